$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

# The existing row 3 (question id) was "1" (same as row 2). It should now be "2",
# since a new row 4 is being added as the new "second question" scenario.
$ws.Range("A3").Value = 2

# Add a new (empty) row 4 below, copying the formatting of row 3, to hold the
# new question for the "carga ePrepago con cuenta de ahorro en cero" scenario.
$ws.Range("A3:O3").Copy() | Out-Null
$ws.Range("A4:O4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Clear any copied values so the new row stays blank, as in the target sheet.
$ws.Range("A4:O4").ClearContents() | Out-Null

# Match the row height used by the other data rows (15.5pt, same as rows 2-3).
$ws.Rows.Item(4).RowHeight = 15.5

# Update the selection to match the new active cell / selected range.
$ws.Range("A4:P4").Select() | Out-Null
